$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 8) for "7号" with count 156,
# matching the style/formatting already used by the row above (row 7).
$ws.Range("A7:B7").Copy() | Out-Null
$ws.Range("A8:B8").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item(8, 1).Value = "7号"
$ws.Cells.Item(8, 2).Value = 156

# Update the selected cell to D6, as recorded in the saved workbook view.
$ws.Range("D6").Select() | Out-Null
